$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "67.272.61"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.874.19"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'469.92"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +9.47%  "
$ws.Range("D6").Value = "'145.28"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +10.46%  "
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.747"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "'0.0000311"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -6.15%  "
$ws.Range("D12").Value = "'43.44"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "4.499.48"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'14.83"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -5.56%  "
$ws.Range("D16").Value = "3.867.53"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "'20.10"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +6.13%  "
$ws.Range("D20").Value = "67.529.72"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "'435.78"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +4.85%  "
$ws.Range("D22").Value = "'14.96"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("E23").Value = "  +5.64%  "
$ws.Range("D24").Value = "'89.32"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("E25").Value = "  +9.79%  "
$ws.Range("D26").Value = "'38.14"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +6.83%  "
$ws.Range("D28").Value = "'9.94"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").Value = "'730.87"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").Value = "'13.88"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  +6.41%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "'44.35"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +12.93%  "
$ws.Range("E35").Value = "  +8.85%  "
$ws.Range("D36").Value = "'58.09"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +4.46%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'5.55"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'0.0485"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.348"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +6.43%  "
$ws.Range("D42").Value = "0.0₃0685"
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'2.55"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +4.49%  "
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  +5.06%  "
$ws.Range("D49").Value = "'2.16"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "'144.48"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +1.70%  "
